$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new product rows (posting multiple products at a time)
$ws.Range("A3").Value = "iPhone 17"
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 120000
$ws.Range("D3").Value = "i7"
$ws.Range("E3").Value = "2TB"

$ws.Range("A4").Value = "One plus 5"
$ws.Range("B4").Value = 2023
$ws.Range("C4").Value = 70000
$ws.Range("D4").Value = "i5"
$ws.Range("E4").Value = "4GB"

# Select one of the newly posted rows
$ws.Range("C11").Select()
